$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing header in B1
$ws.Range("B1").Value = "Category"

# Add a new header cell C1 with the same formatting as B1 (bold, centered, bordered)
$ws.Range("B1").Copy($ws.Range("C1"))
$ws.Range("C1").Value = "Original"

# Fill the new column C with the original category labels for rows 2-7
$ws.Range("C2").Value = "No cat"
$ws.Range("C3").Value = "1 Energy"
$ws.Range("C4").Value = "2 Industrial processes and product use"
$ws.Range("C5").Value = "3 Agriculture"
$ws.Range("C6").Value = "4 Land use, land-use change and forestry (LULUCF)"
$ws.Range("C7").Value = "5 Waste management"
